# Replace "OIE" with "WOAH" across the workbook's textual content, as
# described in the commit message "OIE replaced with WOAH all Excels".
#
# The workbook has two sheets:
#   "Sheet 1"    -> contains narrative text referencing OIE in column E
#   "References" -> contains reference citations referencing OIE in column C
#
# Only whole-word occurrences of "OIE" (including inside the {ref...:OIE}
# placeholders and "OIE-WAHIS"/"OIE-prescribed" compounds) are replaced by
# "WOAH". URLs such as https://wahis.oie.int/ or https://www.oie.int/...
# must remain untouched.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet 1")
$references = $wb.Worksheets.Item("References")

$sheet1.Range("E5").Value = "Based on official disease reports to the WOAH"
$sheet1.Range("E6").Value = "Leishmaniasis is a disease listed in the World Organisation for Animal Health ({ref005:WOAH}) Terrestrial Animal Health Code. The map to the right displays occurrence reported to the {ref001:WOAH-WAHIS} system since 2005."
$sheet1.Range("E7").Value = "As described in the WOAH {ref005:Terrestrial Animal Health Code}, the WOAH early warning system includes immediate notifications and follow-up reports on:"
$sheet1.Range("E14").Value = "Information on stable situations (disease present or absent in a zone or country) is provided by countries through the WOAH monitoring system, which is a different reporting channel. This information is available in a different spatial and temporal scale, which can be browsed on the map independently from the outbreak notification points."
$sheet1.Range("E17").Value = "For more up to date reports, visit the original data source: {ref001:WOAH-WAHIS}."
$sheet1.Range("E72").Value = "WOAH-prescribed tests for international trade include ({ref010:WOAH, Terrestrial Manual}):"
$sheet1.Range("E139").Value = "Geographical distribution data has been kindly provided by the World Organisation of Animal Health (WOAH). {ref001:WOAH-WAHIS} (WOAH World Animal Health Information System) is the original source of these data."

$references.Range("C2").Value = "WOAH-WAHIS (WOAH World Animal Health Information System)"
$references.Range("C5").Value = "WOAH (World Organisation for Animal Health). Terrestrial Animal Health Code 2021. WOAH, Paris, France"
$references.Range("C10").Value = "WOAH (World Organisation for Animal Health), 2021. Leishmaniosis. Chapter 3.1.11. WOAH Terrestrial Manual, Paris, France"
